# Apply the "Updated cryptos list" data refresh to Sheet1.
# For each changed cell: B/C/E are plain text, written directly via .Value.
# D (Price) cells are forced to Text (NumberFormat "@") before the write so that
# values such as "238.00", "0.0692", "1.01" are stored as literal strings -
# matching the source inlineStr cells - instead of being auto-coerced to numbers
# by Excel (which would drop trailing zeros / reparse the dotted thousands format).
# The style is reset to "Normal" right after so the cell keeps its original (default)
# styling, i.e. only the NumberFormat is used transiently for type coercion.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: Bitcoin
$ws.Cells.Item(2,4).NumberFormat = "@"
$ws.Cells.Item(2,4).Value = "35.116.69"
$ws.Cells.Item(2,4).Style = "Normal"
$ws.Cells.Item(2,5).Value = "  +0.41%  "
# Row 3: Ethereum
$ws.Cells.Item(3,4).NumberFormat = "@"
$ws.Cells.Item(3,4).Value = "1.851.61"
$ws.Cells.Item(3,4).Style = "Normal"
$ws.Cells.Item(3,5).Value = "  +1.74%  "
# Row 4: TetherUSD
$ws.Cells.Item(4,5).Value = "  +0.42%  "
# Row 5: BNB
$ws.Cells.Item(5,4).NumberFormat = "@"
$ws.Cells.Item(5,4).Value = "238.00"
$ws.Cells.Item(5,4).Style = "Normal"
$ws.Cells.Item(5,5).Value = "  +3.29%  "
# Row 6: XRP
$ws.Cells.Item(6,5).Value = "  +0.70%  "
# Row 7: USDC
$ws.Cells.Item(7,5).Value = "  +0.40%  "
# Row 8: Solana
$ws.Cells.Item(8,4).NumberFormat = "@"
$ws.Cells.Item(8,4).Value = "42.25"
$ws.Cells.Item(8,4).Style = "Normal"
$ws.Cells.Item(8,5).Value = "  +5.09%  "
# Row 9: Cardano
$ws.Cells.Item(9,5).Value = "  +1.11%  "
# Row 10: Dogecoin
$ws.Cells.Item(10,4).NumberFormat = "@"
$ws.Cells.Item(10,4).Value = "0.0692"
$ws.Cells.Item(10,4).Style = "Normal"
$ws.Cells.Item(10,5).Value = "  +1.38%  "
# Row 11: TRON
$ws.Cells.Item(11,4).NumberFormat = "@"
$ws.Cells.Item(11,4).Value = "0.0990"
$ws.Cells.Item(11,4).Style = "Normal"
$ws.Cells.Item(11,5).Value = "  +0.01%  "
# Row 12: WrappedliquidstakedEther2.0
$ws.Cells.Item(12,4).NumberFormat = "@"
$ws.Cells.Item(12,4).Value = "2.118.54"
$ws.Cells.Item(12,4).Style = "Normal"
$ws.Cells.Item(12,5).Value = "  +1.74%  "
# Row 13: Chainlink
$ws.Cells.Item(13,2).Value = "Chainlink"
$ws.Cells.Item(13,3).Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Cells.Item(13,4).NumberFormat = "@"
$ws.Cells.Item(13,4).Value = "11.40"
$ws.Cells.Item(13,4).Style = "Normal"
$ws.Cells.Item(13,5).Value = "  +0.76%  "
# Row 14: WrappedEther
$ws.Cells.Item(14,2).Value = "WrappedEther"
$ws.Cells.Item(14,3).Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Cells.Item(14,4).NumberFormat = "@"
$ws.Cells.Item(14,4).Value = "1.845.46"
$ws.Cells.Item(14,4).Style = "Normal"
$ws.Cells.Item(14,5).Value = "  +1.39%  "
# Row 15: Polygon
$ws.Cells.Item(15,5).Value = "  +0.90%  "
# Row 16: Polkadot
$ws.Cells.Item(16,5).Value = "  +3.13%  "
# Row 17: WrappedBTC
$ws.Cells.Item(17,4).NumberFormat = "@"
$ws.Cells.Item(17,4).Value = "35.086.59"
$ws.Cells.Item(17,4).Style = "Normal"
# Row 18: Litecoin
$ws.Cells.Item(18,4).NumberFormat = "@"
$ws.Cells.Item(18,4).Value = "70.03"
$ws.Cells.Item(18,4).Style = "Normal"
$ws.Cells.Item(18,5).Value = "  +0.48%  "
# Row 19: ShibaInu
$ws.Cells.Item(19,4).NumberFormat = "@"
$ws.Cells.Item(19,4).Value = "0.0₃0792"
$ws.Cells.Item(19,4).Style = "Normal"
$ws.Cells.Item(19,5).Value = "  +0.96%  "
# Row 20: BitcoinCash
$ws.Cells.Item(20,4).NumberFormat = "@"
$ws.Cells.Item(20,4).Value = "240.53"
$ws.Cells.Item(20,4).Style = "Normal"
$ws.Cells.Item(20,5).Value = "  -0.10%  "
# Row 21: Avalanche
$ws.Cells.Item(21,4).NumberFormat = "@"
$ws.Cells.Item(21,4).Value = "12.15"
$ws.Cells.Item(21,4).Style = "Normal"
$ws.Cells.Item(21,5).Value = "  +0.88%  "
# Row 22: Uniswap
$ws.Cells.Item(22,5).Value = "  +2.01%  "
# Row 23: Dai
$ws.Cells.Item(23,5).Value = "  +0.38%  "
# Row 24: Toncoin
$ws.Cells.Item(24,5).Value = "  -0.26%  "
# Row 25: Monero
$ws.Cells.Item(25,4).NumberFormat = "@"
$ws.Cells.Item(25,4).Value = "169.84"
$ws.Cells.Item(25,4).Style = "Normal"
$ws.Cells.Item(25,5).Value = "  -1.97%  "
# Row 26: Cosmos
$ws.Cells.Item(26,4).NumberFormat = "@"
$ws.Cells.Item(26,4).Value = "8.02"
$ws.Cells.Item(26,4).Style = "Normal"
$ws.Cells.Item(26,5).Value = "  +2.51%  "
# Row 27: PancakeSwap
$ws.Cells.Item(27,4).NumberFormat = "@"
$ws.Cells.Item(27,4).Value = "1.82"
$ws.Cells.Item(27,4).Style = "Normal"
$ws.Cells.Item(27,5).Value = "  +20.83%  "
# Row 28: EthereumClassic
$ws.Cells.Item(28,5).Value = "  +1.46%  "
# Row 29: Stellar
$ws.Cells.Item(29,5).Value = "  +0.19%  "
# Row 30: BinanceUSD
$ws.Cells.Item(30,5).Value = "  +0.45%  "
# Row 31: Hedera
$ws.Cells.Item(31,4).NumberFormat = "@"
$ws.Cells.Item(31,4).Value = "0.0553"
$ws.Cells.Item(31,4).Style = "Normal"
$ws.Cells.Item(31,5).Value = "  +0.76%  "
# Row 32: Filecoin
$ws.Cells.Item(32,4).NumberFormat = "@"
$ws.Cells.Item(32,4).Value = "3.99"
$ws.Cells.Item(32,4).Style = "Normal"
$ws.Cells.Item(32,5).Value = "  -0.32%  "
# Row 33: InternetComputer(DFINITY)
$ws.Cells.Item(33,5).Value = "  +1.46%  "
# Row 34: WEMIXToken
$ws.Cells.Item(34,4).NumberFormat = "@"
$ws.Cells.Item(34,4).Value = "1.72"
$ws.Cells.Item(34,4).Style = "Normal"
$ws.Cells.Item(34,5).Value = "  +26.51%  "
# Row 35: LidoDAOToken
$ws.Cells.Item(35,5).Value = "  +9.64%  "
# Row 36: ImmutableX
$ws.Cells.Item(36,4).NumberFormat = "@"
$ws.Cells.Item(36,4).Value = "0.797"
$ws.Cells.Item(36,4).Style = "Normal"
$ws.Cells.Item(36,5).Value = "  +15.32%  "
# Row 37: TrustWalletToken
$ws.Cells.Item(37,5).Value = "  +3.61%  "
# Row 38: ARBITRUM
$ws.Cells.Item(38,5).Value = "  +9.42%  "
# Row 39: VeChain
$ws.Cells.Item(39,5).Value = "  +4.21%  "
# Row 40: Aave
$ws.Cells.Item(40,4).NumberFormat = "@"
$ws.Cells.Item(40,4).Value = "90.30"
$ws.Cells.Item(40,4).Style = "Normal"
$ws.Cells.Item(40,5).Value = "  -2.69%  "
# Row 41: Maker
$ws.Cells.Item(41,4).NumberFormat = "@"
$ws.Cells.Item(41,4).Value = "1.344.90"
$ws.Cells.Item(41,4).Style = "Normal"
$ws.Cells.Item(41,5).Value = "  +0.46%  "
# Row 42: Gas
$ws.Cells.Item(42,2).Value = "Gas"
$ws.Cells.Item(42,3).Value = "https://coinranking.com/coin/hfw0nnnLtSFc7+gas-gas"
$ws.Cells.Item(42,4).NumberFormat = "@"
$ws.Cells.Item(42,4).Value = "13.34"
$ws.Cells.Item(42,4).Style = "Normal"
$ws.Cells.Item(42,5).Value = "  +57.22%  "
# Row 43: InjectiveProtocol
$ws.Cells.Item(43,2).Value = "InjectiveProtocol"
$ws.Cells.Item(43,3).Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Cells.Item(43,4).NumberFormat = "@"
$ws.Cells.Item(43,4).Value = "14.99"
$ws.Cells.Item(43,4).Style = "Normal"
$ws.Cells.Item(43,5).Value = "  +3.29%  "
# Row 44: RenderToken
$ws.Cells.Item(44,4).NumberFormat = "@"
$ws.Cells.Item(44,4).Value = "2.31"
$ws.Cells.Item(44,4).Style = "Normal"
$ws.Cells.Item(44,5).Value = "  +1.13%  "
# Row 45: HuobiToken
$ws.Cells.Item(45,4).NumberFormat = "@"
$ws.Cells.Item(45,4).Value = "2.45"
$ws.Cells.Item(45,4).Style = "Normal"
$ws.Cells.Item(45,5).Value = "  +1.36%  "
# Row 46: Kaspa
$ws.Cells.Item(46,2).Value = "Kaspa"
$ws.Cells.Item(46,3).Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Cells.Item(46,4).NumberFormat = "@"
$ws.Cells.Item(46,4).Value = "0.0553"
$ws.Cells.Item(46,4).Style = "Normal"
$ws.Cells.Item(46,5).Value = "  +6.17%  "
# Row 47: MXToken
$ws.Cells.Item(47,2).Value = "MXToken"
$ws.Cells.Item(47,3).Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Cells.Item(47,4).NumberFormat = "@"
$ws.Cells.Item(47,4).Value = "2.73"
$ws.Cells.Item(47,4).Style = "Normal"
$ws.Cells.Item(47,5).Value = "  -0.90%  "
# Row 48: FraxShare
$ws.Cells.Item(48,4).NumberFormat = "@"
$ws.Cells.Item(48,4).Value = "6.47"
$ws.Cells.Item(48,4).Style = "Normal"
$ws.Cells.Item(48,5).Value = "  +4.29%  "
# Row 49: RocketPoolETH
$ws.Cells.Item(49,4).NumberFormat = "@"
$ws.Cells.Item(49,4).Value = "2.032.85"
$ws.Cells.Item(49,4).Style = "Normal"
$ws.Cells.Item(49,5).Value = "  +1.72%  "
# Row 50: PaxDollar
$ws.Cells.Item(50,2).Value = "PaxDollar"
$ws.Cells.Item(50,3).Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Cells.Item(50,4).NumberFormat = "@"
$ws.Cells.Item(50,4).Value = "1.01"
$ws.Cells.Item(50,4).Style = "Normal"
$ws.Cells.Item(50,5).Value = "  +0.48%  "
# Row 51: Cronos
$ws.Cells.Item(51,2).Value = "Cronos"
$ws.Cells.Item(51,3).Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Cells.Item(51,4).NumberFormat = "@"
$ws.Cells.Item(51,4).Value = "0.0674"
$ws.Cells.Item(51,4).Style = "Normal"
$ws.Cells.Item(51,5).Value = "  +1.51%  "
